# Add a new "WekWikGeneIdSource" sheet, positioned right after "AddGeneIdSource"
# (and thus right before "BenchlingUrlSource"), mirroring the structure of the
# existing IGEMSource / RepositoryIdSource sheets.
# See: https://github.com/manulera/OpenCloning_backend/issues/263

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("AddGeneIdSource")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "WekWikGeneIdSource"

$headers = @("sequence_file_url", "repository_id", "repository_name", "input", "output", "type", "output_name", "id")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$validation = $newSheet.Range("C2:C1048576").Validation
$validation.Add(3, 1, 1, """addgene,genbank,benchling,snapgene,euroscarf,igem""")
$validation.IgnoreBlank = $true
$validation.InCellDropdown = $true
$validation.ShowInput = $false
$validation.ShowError = $false

# Keep the originally-active sheet selected (adding/renaming a sheet makes it
# active as a side effect).
$wb.Worksheets.Item("NamedThing").Activate() | Out-Null
$wb.Worksheets.Item("NamedThing").Range("A1").Select() | Out-Null
